$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the list of quarter labels (1987Q4 .. 2024Q4) matching the
# existing Dec-31 year-end dates currently stored in A2:A39.
$labels = @()
for ($y = 1987; $y -le 2024; $y++) {
    $labels += ("{0}Q4" -f $y)
}

# The header cell (A1) already carries the plain "text" style (s=1,
# bold font + border, no date number format). Copy that formatting
# down onto A2:A39 so the reused style index matches the header's,
# rather than minting a brand-new style.
$ws.Range("A1").Copy()
$ws.Range("A2:A39").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Replace each date serial with its corresponding "YYYYQ4" text label.
for ($i = 0; $i -lt $labels.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
